# "1st changes of mifos to finflux"
# Insert a new (blank-header) column into the "Repayment schedule" sheet
# between the existing "In Advance" (M) and "Late" (N) columns, shifting
# Late/Outstanding/Disbursement one column to the right, and make the
# "Repayment schedule" sheet the active/selected sheet (it was previously
# "Transactions").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column at N; everything from N onward (Late, Outstanding,
# Disbursement) shifts right to O, P, Q.
$ws.Columns("N").Insert() | Out-Null

# Give the freshly inserted column roughly the same width as its left
# neighbour ("In Advance", column M).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab and move the selection,
# matching the new tabSelected/selection state captured in the workbook.
$ws.Activate() | Out-Null
$ws.Range("S4").Select() | Out-Null
